$wb = $excel.ActiveWorkbook

# --- Rename existing sheets ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Name = "mementos"
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "missing"

# --- Sheet1 ("mementos") edits ---
# B25 becomes a hyperlink to the cv.html archive URL (same text already in the cell)
$ws1.Hyperlinks.Add($ws1.Range("B25"), "https://web.archive.org/web/19990428113621im_/http://dewey.rug.ac.be/barn/tex/cv.html") | Out-Null
$ws1.Range("B25").Style = "Hyperlink"

# --- Sheet2 ("missing") edits: replace old missing-media rows with a summary table ---
$ws2.Range("A3:B5").ClearContents()
$ws2.Columns.Item(2).ColumnWidth = 42.33

# Order matches the shared-string creation order of the source edit
$ws2.Range("A7").Value = "AUDIO"
$ws2.Range("B7").Value = "all missing"

$ws2.Range("A10").Value = "VIDEO"
$ws2.Range("B10").Value = "archived and downloaded but unsupported format"

$ws2.Range("A4").Value = "IMAGES"
$ws2.Range("B4").Value = "none missing"

$ws2.Range("B1").Value = "none missing"

# --- Add Sheet3 ("sputnick") with the music-links memento table ---
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "sputnick"

$ws3.Columns.Item(1).ColumnWidth = 68.5
$ws3.Columns.Item(3).ColumnWidth = 107.83

# Order matches the shared-string creation order of the source edit
$ws3.Range("A1").Value = "original"
$ws3.Range("C1").Value = "memento"

$ws3.Range("A3").Value = "http://www.missouri.edu/~uc489745/music.html"

$ws3.Range("A4").Value = "http://www.yahoo.com/Entertainment/Music/"

$ws3.Range("A5").Value = "http://www.music.indiana.edu/misc/music_resources.html"
$ws3.Range("C5").Value = "https://web.archive.org/web/19981205170619/http://www.music.indiana.edu/misc/music_resources.html"

$ws3.Range("C6").Value = "https://web.archive.org/web/19961120201315/http://syy.oulu.fi/music/"
$ws3.Range("A6").Value = "http://syy.oulu.fi/music.html"

$ws3.Range("A7").Value = "http://www.pathfinder.com/@@2k6FXQAAAAAAgAHU/vibe/mmm/music.html"

$ws3.Range("A8").Value = "http://orpheus.ucsd.edu/webmaster/harmony.html"

$ws3.Range("A9").Value = "http://www.leeds.ac.uk/music.html"
$ws3.Range("C9").Value = "http://web.archive.org/web/19961019180127/http://www.leeds.ac.uk/music.html"

$ws3.Range("A10").Value = "http://harmony-central.mit.edu/"
$ws3.Range("C10").Value = "http://web.archive.org/web/19961105182741/http://harmony-central.com/"
$ws3.Range("D10").Value = "became commercial but shows MIT origins"

$ws3.Range("A11").Value = "http://datura.cerl.uiuc.edu/netstuff/sigsoundLinks.html"
$ws3.Range("C11").Value = "https://web.archive.org/web/19971210111447/http://datura.cerl.uiuc.edu/netstuff/sigsoundLinks.html"

$ws3.Range("A12").Value = "http://datura.cerl.uiuc.edu/schools/courses.html"
$ws3.Range("C12").Value = "https://web.archive.org/web/19971210110432/http://datura.cerl.uiuc.edu/schools/courses.html"

$ws3.Range("A13").Value = "http://american.recordings.com/wwwofmusic/index.html"
$ws3.Range("C13").Value = "https://web.archive.org/web/19961227203847/http://ubl.com/"
$ws3.Range("D13").Value = "became commercial; not sure whether original intent is reflected"

$ws3.Hyperlinks.Add($ws3.Range("A13"), "http://american.recordings.com/wwwofmusic/index.html") | Out-Null
$ws3.Range("A13").Style = "Hyperlink"

# --- Selections / active sheet (order matters: last one activated wins the tab) ---
$ws1.Activate()
$ws1.Range("E14").Select() | Out-Null

$ws2.Activate()
$ws2.Range("B15").Select() | Out-Null

$ws3.Activate()
$ws3.Range("A14").Select() | Out-Null
